# Added aggregated count of services for categories, duration in program.
#
# The worksheet's header row (row 1) previously ended with:
#   ... | Biomedical - Other | Social Protection - Other | Exited | Date Exited
# Six new header columns are inserted right before "Exited" / "Date Exited"
# (which shift from CE1:CF1 to CK1:CL1):
#   Behavioral | Bio-Medical | Post-GBV Care | Social Protection |
#   Other Interventions | duration_in_dreams_program

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 6 blank columns starting at CE, pushing the existing "Exited" /
# "Date Exited" columns (and everything after) to the right.
$ws.Range("CE1:CJ1").EntireColumn.Insert()

# Populate the newly inserted header cells.
$ws.Range("CE1").Value = "Behavioral"
$ws.Range("CF1").Value = "Bio-Medical"
$ws.Range("CG1").Value = "Post-GBV Care"
$ws.Range("CH1").Value = "Social Protection"
$ws.Range("CI1").Value = "Other Interventions"
$ws.Range("CJ1").Value = "duration_in_dreams_program"
